# "wrapping up test file audit"
#
# The author finished auditing this workbook:
#   - removed a stray leftover row ("Sheet" / 3 / 4) from the
#     optimization_parameters sheet (row 16), which had been left behind
#     from earlier testing and didn't belong with the real parameter data
#   - left the newly-selected whole row selected on that sheet
#   - browsed through network_weights (landing on E11)
#   - finished up on the threshold_b tab, which is the sheet left active
#     when the file was saved

$wb = $excel.ActiveWorkbook

# Browse to network_weights and leave the selection on E11.
$wsWeights = $wb.Worksheets.Item("network_weights")
$wsWeights.Activate()
$wsWeights.Range("E11").Select()

# Clean up the stray audit row (A16:"Sheet", B16:3, C16:4) from
# optimization_parameters; everything below shifts up one row.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Activate()
$wsParams.Rows("16:16").Delete()
$wsParams.Range("A16:XFD16").Select()

# Finish on threshold_b -- it's the active tab when the workbook is saved.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
